$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 93. This shifts the existing rows 93..122
# down to 94..123 (carrying all their current values/styles with them).
$ws.Rows("93:93").Insert()

# Populate the newly inserted row 93 with the new record.
$ws.Cells.Item(93, 1).Value = 4
$ws.Cells.Item(93, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(93, 3).Value = "Los Lagos"
$ws.Cells.Item(93, 4).NumberFormat = $ws.Cells.Item(94, 4).NumberFormat
$ws.Cells.Item(93, 4).Value = 44754
$ws.Cells.Item(93, 5).Value = 10
$ws.Cells.Item(93, 6).Value = 100112022
$ws.Cells.Item(93, 7).Value = "Arveja Verde"
$ws.Cells.Item(93, 8).Value = "Perfection"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 70
$ws.Cells.Item(93, 11).Value = 43000
$ws.Cells.Item(93, 12).Value = 43000
$ws.Cells.Item(93, 13).Value = 43000
$ws.Cells.Item(93, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(93, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(93, 16).Value = 1720
$ws.Cells.Item(93, 17).Value = 25
$ws.Cells.Item(93, 18).Value = "Hortaliza"
